$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the confidential disclosure text date from 2021-05-27 to 2021-05-28
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-28 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-8
$ws.Range("D2").Value = 0.4991273560156367
$ws.Range("E2").Value = 0.0009350788137856281

$ws.Range("D3").Value = 0.2448214347274801
$ws.Range("E3").Value = 0.002763636363636301

$ws.Range("D4").Value = 0.09556161942417998
$ws.Range("E4").Value = 0.0011238761238761

$ws.Range("D5").Value = 0.1028723568282355
$ws.Range("E5").Value = 0.001841281531946182

$ws.Range("D6").Value = 0.03033891620912997
$ws.Range("E6").Value = -0.003277767372167117

$ws.Range("D7").Value = 0.02727831679533761
$ws.Range("E7").Value = -0.002520623281393131

$ws.Range("E8").Value = 0.001271934958758214
